# manual_adjust post fix sum transfer
# Rename row label in A2 and refresh the recomputed per-country values in row 2 (B2:EO2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "mean_transfer_over_gdp"

$ws.Range("B2").Value = 0.0413938306792532
$ws.Range("C2").Value = 0.0141816551975709
$ws.Range("D2").Value = 0.00946516095329444
$ws.Range("E2").Value = 0.0000951480991964588
$ws.Range("F2").Value = 0.00570157453184357
$ws.Range("G2").Value = -0.00558446480532269
$ws.Range("H2").Value = 0.0466764060485216
$ws.Range("I2").Value = -0.0117502348337724
$ws.Range("J2").Value = 0.0110454972958553
$ws.Range("K2").Value = 0.0201651903787048
$ws.Range("L2").Value = 0.00909981356802385
$ws.Range("M2").Value = -0.000614566868231799
$ws.Range("N2").Value = 0.000596282677563316
$ws.Range("O2").Value = -0.0090000876383852
$ws.Range("P2").Value = 0.00975849017381169
$ws.Range("Q2").Value = 0.00201405372162403
$ws.Range("R2").Value = 0.00707880875164765
$ws.Range("S2").Value = 0.00440608485326436
$ws.Range("T2").Value = 0.0000697596300239161
$ws.Range("U2").Value = 0.00279507500689928
$ws.Range("V2").Value = 0.0540020919258936
$ws.Range("W2").Value = -0.00938714470721567
$ws.Range("X2").Value = -0.000906456666114651
$ws.Range("Y2").Value = -0.00695717879217674
$ws.Range("Z2").Value = 0.00318292080179674
$ws.Range("AA2").Value = 0.0122116098361174
$ws.Range("AB2").Value = 0.0244086955570301
$ws.Range("AC2").Value = 0.0115402027407016
$ws.Range("AD2").Value = 0.00561471468871215
$ws.Range("AE2").Value = 0.019851765516985
$ws.Range("AF2").Value = 0.010803810441459
$ws.Range("AG2").Value = 0.00403852327272275
$ws.Range("AH2").Value = 0.00925337714862955
$ws.Range("AI2").Value = 0.00108322795466306
$ws.Range("AJ2").Value = -0.00500658095829559
$ws.Range("AK2").Value = -0.00587626998158726
$ws.Range("AL2").Value = 0.0154741353681295
$ws.Range("AM2").Value = -0.00451588879885965
$ws.Range("AN2").Value = 0.00225786505456608
$ws.Range("AO2").Value = -0.000600497116212502
$ws.Range("AP2").Value = 0.00154689264000408
$ws.Range("AQ2").Value = -0.000233546381607169
$ws.Range("AR2").Value = 0.126207069066058
$ws.Range("AS2").Value = -0.00401873659003747
$ws.Range("AT2").Value = -0.0078412996344824
$ws.Range("AU2").Value = 0.017336762032432
$ws.Range("AV2").Value = -0.00635801641070753
$ws.Range("AW2").Value = -0.00406078344217508
$ws.Range("AX2").Value = 0.00352505633088435
$ws.Range("AY2").Value = -0.00524855333424494
$ws.Range("AZ2").Value = 0.0000289666060314453
$ws.Range("BA2").Value = 0.0118990747943195
$ws.Range("BB2").Value = 0.0141833634476714
$ws.Range("BC2").Value = 0.0206760931891854
$ws.Range("BD2").Value = 0.0193347431369358
$ws.Range("BE2").Value = 0.000818420431837212
$ws.Range("BF2").Value = -0.00491683176873893
$ws.Range("BG2").Value = 0.00519281975663479
$ws.Range("BH2").Value = 0.00318990704348918
$ws.Range("BI2").Value = 0.0106347377438979
$ws.Range("BJ2").Value = 0.00266474871843922
$ws.Range("BK2").Value = 0.0271336846630891
$ws.Range("BL2").Value = -0.00603801781059794
$ws.Range("BM2").Value = -0.000268201758402401
$ws.Range("BN2").Value = 0.0030050304162045
$ws.Range("BO2").Value = -0.00282755398221118
$ws.Range("BP2").Value = -0.00950535657251373
$ws.Range("BQ2").Value = -0.00543362439791557
$ws.Range("BR2").Value = -0.00670007859783028
$ws.Range("BS2").Value = -0.00477957394433243
$ws.Range("BT2").Value = 0.00897559131079077
$ws.Range("BU2").Value = -0.000962947433063889
$ws.Range("BV2").Value = -0.0024773543506846
$ws.Range("BW2").Value = 0.015064405996304
$ws.Range("BX2").Value = -0.00102448425189035
$ws.Range("BY2").Value = -0.00929553248775824
$ws.Range("BZ2").Value = -0.00984292562693678
$ws.Range("CA2").Value = 0.0362906630147998
$ws.Range("CB2").Value = -0.00193661000946414
$ws.Range("CC2").Value = 0.00176802813629976
$ws.Range("CD2").Value = 0.0175641219813685
$ws.Range("CE2").Value = -0.00620934815514835
$ws.Range("CF2").Value = -0.00518009841146406
$ws.Range("CG2").Value = -0.00886440422467589
$ws.Range("CH2").Value = 0.00123563251979544
$ws.Range("CI2").Value = -0.0041397522227783
$ws.Range("CJ2").Value = 0.0265627884023923
$ws.Range("CK2").Value = 0.000556142824684237
$ws.Range("CL2").Value = 0.000373370310322225
$ws.Range("CM2").Value = 0.0218207436079144
$ws.Range("CN2").Value = -0.0127230324997295
$ws.Range("CO2").Value = 0.0190793745517016
$ws.Range("CP2").Value = -0.0470765908001615
$ws.Range("CQ2").Value = 0.024122335582292
$ws.Range("CR2").Value = 0.00820554873850347
$ws.Range("CS2").Value = 0.00411233350432648
$ws.Range("CT2").Value = 0.0444877632232503
$ws.Range("CU2").Value = -0.00308437358337995
$ws.Range("CV2").Value = 0.00248851482897805
$ws.Range("CW2").Value = 0.0337354487010322
$ws.Range("CX2").Value = 0.00985667921352351
$ws.Range("CY2").Value = 0.00612261163841476
$ws.Range("CZ2").Value = -0.00586034852731582
$ws.Range("DA2").Value = -0.00445129181313277
$ws.Range("DB2").Value = 0.0148199127783998
$ws.Range("DC2").Value = 0.00493873512489233
$ws.Range("DD2").Value = 0.00116064989934569
$ws.Range("DE2").Value = 0.00304186596491245
$ws.Range("DF2").Value = 0.00388963246416255
$ws.Range("DG2").Value = 0.00837324206125399
$ws.Range("DH2").Value = -0.00268143860403829
$ws.Range("DI2").Value = -0.00441264506610488
$ws.Range("DJ2").Value = 0.00168510515869718
$ws.Range("DK2").Value = 0.00144328622335737
$ws.Range("DL2").Value = 0.0253838478729565
$ws.Range("DM2").Value = 0.00549058552739496
$ws.Range("DN2").Value = 0.0108974049950124
$ws.Range("DO2").Value = -0.00867154915547672
$ws.Range("DP2").Value = 0.035386186290345
$ws.Range("DQ2").Value = 0.0133581836310892
$ws.Range("DR2").Value = -0.00397936067700838
$ws.Range("DS2").Value = -0.00153202620773999
$ws.Range("DT2").Value = -0.00189172803418685
$ws.Range("DU2").Value = -0.00369209911200715
$ws.Range("DV2").Value = -0.0039116880717559
$ws.Range("DW2").Value = 0.00799742707693966
$ws.Range("DX2").Value = 0.00227001630721897
$ws.Range("DY2").Value = 0.0298743651470338
$ws.Range("DZ2").Value = 0.00887485960809664
$ws.Range("EA2").Value = -0.000735188787234849
$ws.Range("EB2").Value = -0.0152771608710901
$ws.Range("EC2").Value = 0.0142548253190453
$ws.Range("ED2").Value = 0.00383363549703493
$ws.Range("EE2").Value = -0.00197227514491281
$ws.Range("EF2").Value = 0.0169875835668753
$ws.Range("EG2").Value = 0.0247565183559394
$ws.Range("EH2").Value = -0.00687259229937567
$ws.Range("EI2").Value = 0.00223217824763678
$ws.Range("EJ2").Value = -0.00242635687211382
$ws.Range("EK2").Value = 0.00274728260881245
$ws.Range("EL2").Value = 0.0417343023927944
$ws.Range("EM2").Value = -0.00145577258648795
$ws.Range("EN2").Value = 0.0104602450626872
$ws.Range("EO2").Value = 0.0181627895729295
